$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The left-hand ("sc" anchor) results table used to extend through row 9;
# it now only has 5 name rows (rows 3-7), so rows 8 and 9's old A:H entries
# (name/anchor score/.../normal columns) need to be removed. The J:Q
# ("positive" anchor) table still uses rows 8 and 9, so only A:H is cleared.
$ws.Range("A8:H9").Clear()

$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.7647058823529411
$ws.Range("C3").Value = 26
$ws.Range("D3").Value = 26
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 8
$ws.Range("J3").Value = "best"
$ws.Range("K3").Value = 0.9152542372881356
$ws.Range("L3").Value = 54
$ws.Range("M3").Value = 54
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 5
$ws.Range("A4").Value = "fraud"
$ws.Range("B4").Value = 0.6944444444444444
$ws.Range("C4").Value = 25
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 11
$ws.Range("J4").Value = "love"
$ws.Range("K4").Value = 0.9130434782608695
$ws.Range("L4").Value = 42
$ws.Range("M4").Value = 42
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 4
$ws.Range("A5").Value = "crisis"
$ws.Range("B5").Value = 0.589041095890411
$ws.Range("C5").Value = 172
$ws.Range("D5").Value = 172
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 120
$ws.Range("J5").Value = "great"
$ws.Range("K5").Value = 0.8660714285714286
$ws.Range("L5").Value = 97
$ws.Range("M5").Value = 97
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 15
$ws.Range("A6").Value = "panic"
$ws.Range("B6").Value = 0.2461240310077519
$ws.Range("C6").Value = 127
$ws.Range("D6").Value = 127
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 389
$ws.Range("J6").Value = "interesting"
$ws.Range("K6").Value = 0.8484848484848485
$ws.Range("L6").Value = 28
$ws.Range("M6").Value = 28
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 5
$ws.Range("A7").Value = "sc"
$ws.Range("B7").Value = 0.2169312169312169
$ws.Range("C7").Value = 41
$ws.Range("D7").Value = 41
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 148
$ws.Range("J7").Value = "positive"
$ws.Range("K7").Value = 0.8448275862068966
$ws.Range("L7").Value = 49
$ws.Range("M7").Value = 49
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 9
$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8170731707317073
$ws.Range("L8").Value = 67
$ws.Range("M8").Value = 67
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 15
$ws.Range("J9").Value = "free"
$ws.Range("K9").Value = 0.7833333333333333
$ws.Range("L9").Value = 94
$ws.Range("M9").Value = 94
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 26
$ws.Range("J10").Value = "thank"
$ws.Range("K10").Value = 0.78125
$ws.Range("L10").Value = 100
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 28
$ws.Range("J11").Value = "support"
$ws.Range("K11").Value = 0.7547169811320755
$ws.Range("L11").Value = 80
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 26
$ws.Range("J12").Value = "special"
$ws.Range("K12").Value = 0.75
$ws.Range("L12").Value = 27
$ws.Range("M12").Value = 27
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 9
$ws.Range("J13").Value = "won"
$ws.Range("K13").Value = 0.7435897435897436
$ws.Range("L13").Value = 29
$ws.Range("M13").Value = 29
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 10
$ws.Range("J14").Value = "safe"
$ws.Range("K14").Value = 0.7323943661971831
$ws.Range("L14").Value = 104
$ws.Range("M14").Value = 104
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 38
$ws.Range("J15").Value = "good"
$ws.Range("K15").Value = 0.725
$ws.Range("L15").Value = 116
$ws.Range("M15").Value = 116
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 44
$ws.Range("J16").Value = "safety"
$ws.Range("K16").Value = 0.6666666666666666
$ws.Range("L16").Value = 34
$ws.Range("M16").Value = 34
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 17
$ws.Range("J17").Value = "heroes"
$ws.Range("K17").Value = 0.6595744680851063
$ws.Range("L17").Value = 31
$ws.Range("M17").Value = 31
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 16
$ws.Range("J18").Value = "well"
$ws.Range("K18").Value = 0.6382978723404256
$ws.Range("L18").Value = 60
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 34
$ws.Range("J19").Value = "better"
$ws.Range("K19").Value = 0.6031746031746031
$ws.Range("L19").Value = 38
$ws.Range("M19").Value = 38
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 25
$ws.Range("J20").Value = "relief"
$ws.Range("K20").Value = 0.6
$ws.Range("L20").Value = 30
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 20
$ws.Range("J21").Value = "fresh"
$ws.Range("K21").Value = 0.5833333333333334
$ws.Range("L21").Value = 28
$ws.Range("M21").Value = 28
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 20
$ws.Range("J22").Value = "hand"
$ws.Range("K22").Value = 0.5248041775456919
$ws.Range("L22").Value = 201
$ws.Range("M22").Value = 201
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 182
$ws.Range("J23").Value = "help"
$ws.Range("K23").Value = 0.4779661016949153
$ws.Range("L23").Value = 141
$ws.Range("M23").Value = 141
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 154
$ws.Range("J24").Value = "like"
$ws.Range("K24").Value = 0.4705882352941176
$ws.Range("L24").Value = 160
$ws.Range("M24").Value = 160
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 180
$ws.Range("J25").Value = "care"
$ws.Range("K25").Value = 0.4382022471910113
$ws.Range("L25").Value = 39
$ws.Range("M25").Value = 39
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 50
$ws.Range("J26").Value = "protect"
$ws.Range("K26").Value = 0.4246575342465753
$ws.Range("L26").Value = 31
$ws.Range("M26").Value = 31
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 42
$ws.Range("J27").Value = "increase"
$ws.Range("K27").Value = 0.3461538461538461
$ws.Range("L27").Value = 27
$ws.Range("M27").Value = 27
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 51
$ws.Range("J28").Value = "please"
$ws.Range("K28").Value = 0.3347280334728033
$ws.Range("L28").Value = 80
$ws.Range("M28").Value = 80
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 159
$ws.Range("J29").Value = "store"
$ws.Range("K29").Value = 0.03243847874720358
$ws.Range("L29").Value = 29
$ws.Range("M29").Value = 29
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 865
$ws.Range("J30").Value = "you"
$ws.Range("K30").Value = 0.03
$ws.Range("L30").Value = 36
$ws.Range("M30").Value = 36
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = $false
$ws.Range("Q30").Value = 1164
$ws.Range("J31").Value = "!"
$ws.Range("K31").Value = 0.02836879432624113
$ws.Range("L31").Value = 28
$ws.Range("M31").Value = 30
$ws.Range("N31").Value = 0.93
$ws.Range("O31").Value = 0.06999999999999995
$ws.Range("P31").Value = $true
$ws.Range("Q31").Value = 959
$ws.Range("J32").Value = ","
$ws.Range("K32").Value = 0.01556739041376485
$ws.Range("L32").Value = 38
$ws.Range("M32").Value = 42
$ws.Range("N32").Value = 0.9
$ws.Range("O32").Value = 0.09999999999999998
$ws.Range("P32").Value = $true
$ws.Range("Q32").Value = 2403
$ws.Range("J33").Formula = '=TEXT(19,"0")'
$ws.Range("J33").Copy()
$ws.Range("J33").PasteSpecial(-4163)
$ws.Range("K33").Value = 0.01449953227315248
$ws.Range("L33").Value = 31
$ws.Range("M33").Value = 36
$ws.Range("N33").Value = 0.86
$ws.Range("O33").Value = 0.14
$ws.Range("P33").Value = $true
$ws.Range("Q33").Value = 2107
$ws.Range("J34").Value = "."
$ws.Range("K34").Value = 0.01422275641025641
$ws.Range("L34").Value = 71
$ws.Range("M34").Value = 75
$ws.Range("N34").Value = 0.95
$ws.Range("O34").Value = 0.05000000000000004
$ws.Range("P34").Value = $true
$ws.Range("Q34").Value = 4921
$ws.Range("J35").Value = "to"
$ws.Range("K35").Value = 0.01410730804810361
$ws.Range("L35").Value = 61
$ws.Range("M35").Value = 65
$ws.Range("N35").Value = 0.9399999999999999
$ws.Range("O35").Value = 0.06000000000000005
$ws.Range("P35").Value = $true
$ws.Range("Q35").Value = 4263
$ws.Range("J36").Value = "and"
$ws.Range("K36").Value = 0.01386286998875983
$ws.Range("L36").Value = 37
$ws.Range("M36").Value = 41
$ws.Range("N36").Value = 0.9
$ws.Range("O36").Value = 0.09999999999999998
$ws.Range("P36").Value = $true
$ws.Range("Q36").Value = 2632
$ws.Range("J37").Value = "a"
$ws.Range("K37").Value = 0.01258680555555556
$ws.Range("L37").Value = 29
$ws.Range("M37").Value = 34
$ws.Range("N37").Value = 0.85
$ws.Range("O37").Value = 0.15
$ws.Range("P37").Value = $true
$ws.Range("Q37").Value = 2275
$ws.Range("J38").Value = "co"
$ws.Range("K38").Value = 0.01127214170692432
$ws.Range("L38").Value = 35
$ws.Range("M38").Value = 37
$ws.Range("N38").Value = 0.95
$ws.Range("O38").Value = 0.05000000000000004
$ws.Range("P38").Value = $true
$ws.Range("Q38").Value = 3070
$ws.Range("J39").Value = "the"
$ws.Range("K39").Value = 0.01007361487795428
$ws.Range("L39").Value = 52
$ws.Range("M39").Value = 55
$ws.Range("N39").Value = 0.95
$ws.Range("O39").Value = 0.05000000000000004
$ws.Range("P39").Value = $true
$ws.Range("Q39").Value = 5110
